$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = "44.480.79"
$ws.Range("E2").Value = "  +3.71%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = "2.280.65"
$ws.Range("E3").Value = "  +2.52%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.10"
$ws.Range("E5").Value = "  +1.23%  "

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.14"
$ws.Range("E6").Value = "  +5.89%  "

# Row 7: XRP -> XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.593"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8: USDC -> USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9: Cardano -> Cardano
$ws.Range("E9").Value = "  +1.66%  "

# Row 10: Avalanche -> Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.82"
$ws.Range("E10").Value = "  +3.60%  "

# Row 11: Dogecoin -> Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +1.89%  "

# Row 12: Polkadot -> Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.91"
$ws.Range("E12").Value = "  +1.31%  "

# Row 13: TRON -> TRON
$ws.Range("E13").Value = "  +0.89%  "

# Row 14: Polygon -> Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.887"
$ws.Range("E14").Value = "  +2.41%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.632.42"
$ws.Range("E15").Value = "  +2.81%  "

# Row 16: Chainlink -> Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.65"
$ws.Range("E16").Value = "  +2.45%  "

# Row 17: WrappedEther -> WrappedEther
$ws.Range("D17").Value = "2.283.97"
$ws.Range("E17").Value = "  +2.39%  "

# Row 18: WrappedBTC -> WrappedBTC
$ws.Range("D18").Value = "44.365.70"
$ws.Range("E18").Value = "  +3.58%  "

# Row 19: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.11"
$ws.Range("E19").Value = "  -5.86%  "

# Row 20: ShibaInu -> ShibaInu
$ws.Range("E20").Value = "  +4.22%  "

# Row 21: Uniswap -> Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.54"
$ws.Range("E21").Value = "  +1.18%  "

# Row 22: Litecoin -> Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.49"
$ws.Range("E22").Value = "  +1.35%  "

# Row 23: PancakeSwap -> PancakeSwap
$ws.Range("E23").Value = "  +1.75%  "

# Row 24: BitcoinCash -> BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.62"
$ws.Range("E24").Value = "  +1.20%  "

# Row 25: ImmutableX -> ImmutableX
$ws.Range("E25").Value = "  +3.10%  "

# Row 26: Dai -> Dai
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "

# Row 27: Cosmos -> Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("E27").Value = "  +1.24%  "

# Row 28: InjectiveProtocol -> Toncoin
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  +0.38%  "

# Row 29: Toncoin -> InjectiveProtocol
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.41"
$ws.Range("E29").Value = "  +12.18%  "

# Row 30: Filecoin -> Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.53"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31: Monero -> Monero
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.64"
$ws.Range("E31").Value = "  +4.20%  "

# Row 32: EthereumClassic -> EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.65"
$ws.Range("E32").Value = "  +0.78%  "

# Row 33: Hedera -> Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0885"
$ws.Range("E33").Value = "  -2.88%  "

# Row 34: WEMIXToken -> WEMIXToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.75"
$ws.Range("E34").Value = "  -1.04%  "

# Row 35: ARBITRUM -> ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +4.61%  "

# Row 36: Kaspa -> Kaspa
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  +10.24%  "

# Row 37: LidoDAOToken -> LidoDAOToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.20"
$ws.Range("E37").Value = "  +0.39%  "

# Row 38: Stellar -> Stellar
$ws.Range("E38").Value = "  -0.56%  "

# Row 39: NEARProtocol -> NEARProtocol
$ws.Range("E39").Value = "  -0.15%  "

# Row 40: RenderToken -> RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.47"
$ws.Range("E40").Value = "  -0.16%  "

# Row 41: VeChain -> Celestia
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.58"
$ws.Range("E41").Value = "  +24.63%  "

# Row 42: Celestia -> VeChain
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0329"
$ws.Range("E42").Value = "  +0.76%  "

# Row 43: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("E43").Value = "  +0.23%  "

# Row 44: Maker -> Maker
$ws.Range("D44").Value = "1.772.97"
$ws.Range("E44").Value = "  -9.16%  "

# Row 45: Algorand -> BitcoinSV
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "87.06"
$ws.Range("E45").Value = "  -2.16%  "

# Row 46: BitcoinSV -> Algorand
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.208"
$ws.Range("E46").Value = "  -0.02%  "

# Row 47: THORChain -> THORChain
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.51"
$ws.Range("E47").Value = "  +2.25%  "

# Row 48: MultiversX -> MultiversX
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "60.39"
$ws.Range("E48").Value = "  -0.38%  "

# Row 49: Stacks -> ordi
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.80"
$ws.Range("E49").Value = "  -2.05%  "

# Row 50: ordi -> Stacks
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.71"
$ws.Range("E50").Value = "  +4.53%  "

# Row 51: FraxShare -> Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.55"
$ws.Range("E51").Value = "  +1.45%  "
